# Update column C (dates) for rows 2-57 from 45190 (2023-09-21) to 45192 (2023-09-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 57; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
